# Weekly price update: a new week's record is published for
# "Terminal La Palmera de La Serena - Jengibre". The new record is
# inserted as row 71 (most recent date first), which pushes every
# existing record from row 71 downward by one row (71->72, ..., 188->189).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 71, shifting rows 71:188 down to 72:189
$ws.Rows("71:71").Insert()

# Populate the newly inserted row 71 with the new week's data
$ws.Cells.Item(71, 1).Value  = 8
$ws.Cells.Item(71, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(71, 3).Value  = "Coquimbo"
$ws.Cells.Item(71, 4).Value  = 45210
$ws.Cells.Item(71, 5).Value  = 4
$ws.Cells.Item(71, 6).Value  = 100114007
$ws.Cells.Item(71, 7).Value  = "Jengibre"
$ws.Cells.Item(71, 8).Value  = "Sin especificar"
$ws.Cells.Item(71, 9).Value  = "Primera"
$ws.Cells.Item(71, 10).Value = 460
$ws.Cells.Item(71, 11).Value = 21000
$ws.Cells.Item(71, 12).Value = 22000
$ws.Cells.Item(71, 13).Value = 21500
$ws.Cells.Item(71, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(71, 15).Value = "Perú"
$ws.Cells.Item(71, 16).Value = 1654
$ws.Cells.Item(71, 17).Value = 13
$ws.Cells.Item(71, 18).Value = "Hortaliza"
